$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain text by pre-formatting as Text.
$ws.Range("D2").Value = "35.312.73"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.912.23"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +8.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "255.74"
$ws.Range("E6").Value = "  +3.74%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.18"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.366"
$ws.Range("E9").Value = "  +5.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.29"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0769"
$ws.Range("E11").Value = "  +6.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0988"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.15"
$ws.Range("E13").Value = "  +6.66%  "
$ws.Range("D14").Value = "2.189.38"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.740"
$ws.Range("E15").Value = "  +5.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.02"
$ws.Range("E16").Value = "  +4.09%  "
$ws.Range("D17").Value = "1.913.54"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "35.294.74"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.13"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("E20").Value = "  +3.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "245.94"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.16"
$ws.Range("E22").Value = "  +5.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.17"
$ws.Range("E23").Value = "  +7.13%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  +7.70%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.87"
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("D31").Value = "4.127.54"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.69"
$ws.Range("E32").Value = "  +26.66%  "
$ws.Range("E33").Value = "  +5.33%  "
$ws.Range("E34").Value = "  +15.46%  "
$ws.Range("E35").Value = "  +4.80%  "
$ws.Range("E36").Value = "  +4.47%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.921"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "100.58"
$ws.Range("E40").Value = "  +11.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0221"
$ws.Range("E41").Value = "  +6.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.10"
$ws.Range("E42").Value = "  +5.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.14"
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("D46").Value = "1.344.03"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.74"
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.10"
$ws.Range("E50").Value = "  -7.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0758"
$ws.Range("E51").Value = "  +7.24%  "
